# Presentation - Source Images.pptx : "Add files via upload" edit
#
# 1. Bumps the cached datetimeFigureOut field text (3/11/25 -> 3/14/25) on every
#    slide layout + the slide master.
# 2. Re-flows / re-positions + re-words the title and subtitle text boxes on
#    slide 1 ("Forecasting ... " headline and the italic "AI Forecasts ..."
#    strap-line), shrinking the headline font 44pt -> 42pt and nudging both
#    boxes up/left a touch.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder field text on every slide layout + the slide master.
# ---------------------------------------------------------------------------
$targets = New-Object System.Collections.ArrayList
foreach ($m in $p.Designs) {
    $null = $targets.Add($m.SlideMaster)
    foreach ($l in $m.SlideMaster.CustomLayouts) {
        $null = $targets.Add($l)
    }
}

foreach ($master in $p.SlideMaster) {
    $null = $targets.Add($master)
}

foreach ($container in $targets) {
    foreach ($shp in $container.Shapes) {
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "3/11/25") {
                    $tr.Text = "3/14/25"
                }
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 1 content.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# ---- TextBox 5 (headline: "Forecasting the NFL GOAT QB with Artificial
#      Intelligence") ----
$headline = $s.Shapes.Item(2)

# Re-position (EMU 122949,1947128 -> 43926,1800577; size stays the same).
$headline.Left = 3.458740157480315
$headline.Top = 141.77771763543305

$htr = $headline.TextFrame.TextRange

# Paragraph 1 originally holds two runs: "Forecasting" + " ". Extend the
# first run's text (keeps its bg1-coloured rPr) and drop the now-redundant
# second run.
$run1 = $htr.Characters(1, 11)
$run1.Text = "Forecasting the "
$afterRun1 = $headline.TextFrame.TextRange
$oldSpaceRun = $afterRun1.Characters(17, 1)
$oldSpaceRun.Text = ""

# Paragraph 2 first run: "The GOAT " -> "NFL GOAT QB " (rPr / colour unchanged).
$afterTrim = $headline.TextFrame.TextRange
$goatRun = $afterTrim.Characters(18, 9)
$goatRun.Text = "NFL GOAT QB "

# Shrink every run (and the stray trailing paragraph mark) from 44pt to 42pt.
$headline.TextFrame.TextRange.Font.Size = 42

# ---- TextBox 6 (strap-line: "AI Forecasts How the Tom Brady vs. Patrick
#      Mahomes Debate Will End") ----
$strap = $s.Shapes.Item(3)

# Re-position (EMU 122949,4955823 -> 43926,4678833; size stays the same).
$strap.Left = 3.458740157480315
$strap.Top = 368.4120331440945

$stext = $strap.TextFrame.TextRange
$stext.Text = "AI Forecasts How the Tom Brady vs. Patrick Mahomes Debate Will End"
